$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("M2").Value = 1.07
$ws.Range("N2").Value = 9
$ws.Range("Q2").Value = 1.59
$ws.Range("R2").Value = 2.33
$ws.Range("S2").Value = 2.15
$ws.Range("T2").Value = 1.67
$ws.Range("U2").Value = 3.05
$ws.Range("V2").Value = 1.37
$ws.Range("W2").Value = 4
$ws.Range("X2").Value = 1.22

# Row 3 updates
$ws.Range("Q3").Value = 1.83
$ws.Range("R3").Value = 2.03
$ws.Range("S3").Value = 2.4
$ws.Range("T3").Value = 1.53

# Row 5 updates
$ws.Range("M5").Value = 1.07
$ws.Range("N5").Value = 9
$ws.Range("O5").Value = 1.36
$ws.Range("P5").Value = 3
$ws.Range("S5").Value = 2.15
$ws.Range("T5").Value = 1.67
$ws.Range("W5").Value = 4
$ws.Range("X5").Value = 1.22
